$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("BAEPAbCiPC")

# ---------------------------------------------------------------------------
# About sheet: remove the trailing explanatory paragraph (rows 15-18), which
# also drops the now-unused shared strings describing the U.S. model pass
# through behaviour.
# ---------------------------------------------------------------------------
$wsAbout.Rows("15:18").Delete()

# ---------------------------------------------------------------------------
# BAEPAbCiPC sheet: flip several fuel/energy-carrier flags from 0 to 1.
# ---------------------------------------------------------------------------
$wsData.Range("B3").Value = 1   # hard coal
$wsData.Range("B4").Value = 1   # natural gas
$wsData.Range("B9").Value = 1   # biomass
$wsData.Range("B10").Value = 1  # petroleum gasoline
$wsData.Range("B11").Value = 1  # petroleum diesel
$wsData.Range("B12").Value = 1  # biofuel gasoline
$wsData.Range("B13").Value = 1  # biofuel diesel
$wsData.Range("B14").Value = 1  # jet fuel or kerosene
$wsData.Range("B17").Value = 1  # lignite
$wsData.Range("B18").Value = 1  # crude oil
$wsData.Range("B19").Value = 1  # heavy fuel oil
$wsData.Range("B20").Value = 1  # LPG propane or butane

# ---------------------------------------------------------------------------
# Minor view/selection tweaks captured in the diff.
# ---------------------------------------------------------------------------
$wsData.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$wsData.Range("B17:B20").Select()

$wsAbout.Activate()
$wsAbout.Range("A15:XFD18").Select()
